$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column layout cleanup: split the combined A:B column definition so
#     column A keeps its own 30.71-char width/style instead of sharing
#     the stale min=1,max=2 range with column B. ---
$ws.Columns("A").ColumnWidth = 29.83

# --- Row 10 (Objetivos:) had the wrong value (a leftover professor name);
#     replace it with the real Portuguese objectives text. ---
$ws.Range("B10").Value = "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."
$ws.Range("C10").Value = "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."

# --- Insert a new row 13 to hold the "Docentes responsáveis:" value,
#     which was missing (all rows below it were off by one). ---
$ws.Rows("13:13").Insert()
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "5817650 - Érica Leonor Romão"
$ws.Range("C13").Value = "5817650 - Érica Leonor Romão"
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fix the remaining rows whose value cells were each carrying the
#     content meant for the row below (cascading off-by-one data bug). ---

# Row 14 (was 13): Programa resumido:
$ws.Range("B14").Value = "A definir de acordo com o tópico programado"
$ws.Range("C14").Value = "A definir de acordo com o tópico programado"

# Row 16 (was 15): Programa:
$ws.Range("B16").Value = "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."
$ws.Range("C16").Value = "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."

# Row 19 (was 18): Método:
$ws.Range("B19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa. As avalições podem ser: seminários, trabalhos, projetos ou outra forma de avaliação definida pelo professor. Sendo necessário no mínimo uma avaliação na forma de prova escrita."
$ws.Range("C19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa. As avalições podem ser: seminários, trabalhos, projetos ou outra forma de avaliação definida pelo professor. Sendo necessário no mínimo uma avaliação na forma de prova escrita."

# Row 20 (was 19): Critério:
$ws.Range("B20").Value = "Média ponderada das avaliações (M)."
$ws.Range("C20").Value = "Média ponderada das avaliações (M)."

# Row 21 (was 20): Norma de recuperação:
$ws.Range("B21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"
$ws.Range("C21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"

# Row 22 (was 21): Bibliografia: - now gets its own correct value instead of
# the recovery-rule text that had spilled into it.
$ws.Range("B22").Value = "Livros, artigos ou texto fornecido pelo docente responsável extraídos de livros ou revistas especializadas na área de Meio Ambiente."
$ws.Range("C22").Value = "Livros, artigos ou texto fornecido pelo docente responsável extraídos de livros ou revistas especializadas na área de Meio Ambiente."
